# "Generate Report for Archive"
# The localization status report is regenerated: rows that were previously
# marked "Ready for handoff" are now "In Translation" (reflected on the
# Overview sheet's per-locale columns as well as each locale sheet's
# Status column). The narrower status text means the columns that show it
# no longer need to be as wide, so their column widths shrink to fit the
# new content.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# -- Overview sheet: zh-cn / de-de status columns (E, F) --
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus

# Columns shrink to fit the shorter status text.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# -- zh-cn sheet: Status column (C) --
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# -- de-de sheet: Status column (C) --
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
